$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($year = 2004; $year -le 2024; $year++) {
    $row = $year - 2004 + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "general"
    $cell.Value = "$($year)Q4"
}
